$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.835.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.93%  '

$ws.Range("D3").Value = "'3.340.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.39%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = "'573.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.28%  '

$ws.Range("D6").Value = "'181.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.24%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -1.55%  '

$ws.Range("E9").Value = '  -3.76%  '

$ws.Range("D10").Value = "'6.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.74%  '

$ws.Range("E11").Value = '  -4.47%  '

$ws.Range("D12").Value = "'3.919.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.38%  '

$ws.Range("E13").Value = '  -1.81%  '

$ws.Range("D14").Value = "'27.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.61%  '

$ws.Range("D15").Value = "'66.881.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.89%  '

$ws.Range("E16").Value = '  -2.66%  '

$ws.Range("D17").Value = "'3.326.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.08%  '

$ws.Range("D18").Value = "'438.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.89%  '

$ws.Range("D19").Value = "'13.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.81%  '

$ws.Range("D20").Value = "'5.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.87%  '

$ws.Range("E21").Value = '  -2.74%  '

$ws.Range("D22").Value = "'73.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.72%  '

$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").Value = "'0.519"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.78%  '

$ws.Range("E25").Value = '  -4.37%  '

$ws.Range("E26").Value = '  -0.26%  '

$ws.Range("D27").Value = "'9.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.96%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("E29").Value = '  -1.42%  '

$ws.Range("D30").Value = "'22.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.53%  '

$ws.Range("D31").Value = "'5.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.31%  '

$ws.Range("E32").Value = '  -0.01%  '

$ws.Range("D33").Value = "'6.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.75%  '

$ws.Range("D34").Value = "'1.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.62%  '

$ws.Range("D35").Value = "'162.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.76%  '

$ws.Range("D36").Value = "'1.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.00%  '

$ws.Range("D37").Value = "'27.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.16%  '

$ws.Range("D38").Value = "'1.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.84%  '

$ws.Range("D39").Value = "'2.826.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.33%  '

$ws.Range("E40").Value = '  -3.40%  '

$ws.Range("D41").Value = "'4.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.80%  '

$ws.Range("D42").Value = "'6.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.50%  '

$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").Value = "'0.0672"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.01%  '

$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = "'40.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.62%  '

$ws.Range("D45").Value = "'24.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.14%  '

$ws.Range("D46").Value = "'2.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.14%  '

$ws.Range("D47").Value = "'323.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.69%  '

$ws.Range("D48").Value = "'0.0274"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.25%  '

$ws.Range("D49").Value = "'0.990"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.46%  '

$ws.Range("D50").Value = "'6.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.92%  '

$ws.Range("D51").Value = "'30.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.66%  '
